# Agregar partidos de Liga Peruana 2025 (jornadas del 18-20 de julio)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S")
# Indices (0-based, relative to the 19 data columns) that hold text values
# rather than numbers: Fecha, Local, Visita, Posesion Local, Posesion Visita, Resultado
$textCols = @(0, 1, 2, 16, 17, 18)

# Each inner array: row number followed by the 19 column values (A..S)
$rows = @(
    @(173, '2025-07-18', 'Alianza Atletico', 'Sport Huancayo', 0, 0, 1391870, 5, 4, 0, 0, 0, 0, 0, 0, 0, 0, '57%', '43%', 'E'),
    @(174, '2025-07-19', 'Deportivo Garcilaso', 'Cultural Santa Rosa', 3, 0, 1405332, 5, 9, 2, 0, 0, 0, 0, 0, 3, 0, '44%', '56%', 'L'),
    @(175, '2025-07-19', 'Universitario', 'Comerciantes Unidos', 3, 1, 1391871, 6, 2, 0, 2, 0, 0, 0, 0, 3, 1, '59%', '41%', 'L'),
    @(176, '2025-07-19', 'UTC', 'FBC Melgar', 1, 2, 1391872, 8, 7, 2, 3, 1, 0, 0, 0, 1, 2, '45%', '55%', 'V'),
    @(177, '2025-07-19', 'ADT', 'Cienciano', 1, 0, 1391873, 3, 11, 3, 3, 0, 0, 0, 0, 1, 0, '47%', '53%', 'L'),
    @(178, '2025-07-19', 'Cusco', 'Alianza Lima', 2, 0, 1391874, 7, 3, 5, 2, 0, 0, 0, 0, 2, 0, '62%', '38%', 'L'),
    @(179, '2025-07-20', 'Sporting Cristal', 'Alianza Universidad', 3, 0, 1391875, 4, 3, 4, 2, 0, 0, 0, 0, 3, 0, '56%', '44%', 'L'),
    @(180, '2025-07-20', 'Ayacucho FC', 'Atletico Grau', 1, 2, 1391876, 3, 1, 1, 2, 0, 0, 0, 0, 1, 2, '53%', '47%', 'V'),
    @(181, '2025-07-20', 'Juan Pablo II College', 'Sport Boys', 3, 0, 1391877, 1, 6, 1, 2, 0, 1, 0, 0, 3, 0, '43%', '57%', 'L')
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $colLetter = $cols[$i]
        $val = $r[$i + 1]
        $addr = "$colLetter$rowNum"
        if ($textCols -contains $i) {
            # Force text type so values like "2025-07-18" or "57%" are not
            # auto-converted to dates/numbers, while keeping default styling.
            $ws.Range($addr).NumberFormat = "@"
            $ws.Range($addr).Value = $val
            $ws.Range($addr).Style = "Normal"
        } else {
            $ws.Range($addr).Value = $val
        }
    }
}
